# Update "Name of Algo" - refresh imputed values in result_data_RandomForest.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.458499999999999
$ws.Range("B3").Value = 5.880099999999988
$ws.Range("B5").Value = 4.880300000000003
$ws.Range("A9").Value = -20.46299999999997
$ws.Range("B11").Value = 5.347099999999997
$ws.Range("B12").Value = 5.334299999999998
$ws.Range("A13").Value = -21.94500000000002
$ws.Range("A16").Value = -20.12009999999999
$ws.Range("A18").Value = -22.7147
$ws.Range("A20").Value = -22.07480000000001
$ws.Range("B21").Value = 5.483999999999993
